$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B/C rows for each year group (only columns A:G - the "_累计值" data).
# 2012: rows 3 (B) and 4 (C)
$row3 = $ws.Range("A3:G3").Value2
$row4 = $ws.Range("A4:G4").Value2
$ws.Range("A3:G3").Value2 = $row4
$ws.Range("A4:G4").Value2 = $row3

# 2013: rows 6 (B) and 7 (C)
$row6 = $ws.Range("A6:G6").Value2
$row7 = $ws.Range("A7:G7").Value2
$ws.Range("A6:G6").Value2 = $row7
$ws.Range("A7:G7").Value2 = $row6

# 2014: rows 9 (B) and 10 (C)
$row9 = $ws.Range("A9:G9").Value2
$row10 = $ws.Range("A10:G10").Value2
$ws.Range("A9:G9").Value2 = $row10
$ws.Range("A10:G10").Value2 = $row9

# 2015: rows 12 (B) and 13 (C)
$row12 = $ws.Range("A12:G12").Value2
$row13 = $ws.Range("A13:G13").Value2
$ws.Range("A12:G12").Value2 = $row13
$ws.Range("A13:G13").Value2 = $row12

# Remove the duplicated non-cumulative columns H:M entirely (header + data),
# leaving only A:G ("_累计值" columns).
$ws.Range("H1:M16").EntireColumn.Delete()
